$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like pure numbers need to be forced to stay
# as text (matching the source workbook, where every cell is stored as a
# string) - otherwise Excel will auto-convert them to numeric cells.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2395540825"

$ws.Range("L2").Value = ": HEAD-UP-DISPLAY"

$ws.Range("O2").Value = "587,91"

$ws.Range("P2").Value = "ATC00401281230820195604"

$ws.Range("Q2").NumberFormat = "@"
$ws.Range("Q2").Value = "20190826"

$ws.Range("R2").NumberFormat = "@"
$ws.Range("R2").Value = "3012.93"

$ws.Range("S2").Value = "81,35"

$ws.Range("AF2").NumberFormat = "@"
$ws.Range("AF2").Value = "85122000900"
